# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to each profession sheet
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 301.9355
$ws.Cells.Item(12, 9).Value = 147
$ws.Cells.Item(12, 11).Value = 147
$ws.Cells.Item(12, 13).Value = 23
$ws.Cells.Item(33, 8).Value = 11062.934
$ws.Cells.Item(33, 9).Value = 13715.958
$ws.Cells.Item(33, 10).Value = 450.83334
$ws.Cells.Item(33, 11).Value = 13715.958
$ws.Cells.Item(33, 12).Value = 450.83334
$ws.Cells.Item(33, 13).Value = -13486.958
$ws.Cells.Item(33, 14).Value = -908.83334
$ws.Cells.Item(62, 8).Value = 4831.6665
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 4831.6665
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 4831.6665
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -6079.6665
$ws.Cells.Item(65, 8).Value = 4831.6665
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 4831.6665
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 24158.3325
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -30398.3325
$ws.Cells.Item(81, 8).Value = 76646.336
$ws.Cells.Item(81, 10).Value = 76646.336
$ws.Cells.Item(81, 12).Value = 76646.336
$ws.Cells.Item(81, 14).Value = -78642.336
$ws.Cells.Item(84, 8).Value = 76646.336
$ws.Cells.Item(84, 10).Value = 76646.336
$ws.Cells.Item(84, 12).Value = 229939.008
$ws.Cells.Item(84, 14).Value = -239923.008
$ws.Cells.Item(98, 8).Value = 871.05554
$ws.Cells.Item(98, 9).Value = 880.375
$ws.Cells.Item(98, 11).Value = 880.375
$ws.Cells.Item(98, 13).Value = 617.625
$ws.Cells.Item(122, 8).Value = 871.05554
$ws.Cells.Item(122, 9).Value = 880.375
$ws.Cells.Item(122, 11).Value = 2641.125
$ws.Cells.Item(122, 13).Value = -191.125
$ws.Cells.Item(125, 8).Value = 6262975
$ws.Cells.Item(125, 9).Value = 24500
$ws.Cells.Item(125, 10).Value = 12501450
$ws.Cells.Item(125, 11).Value = 220500
$ws.Cells.Item(125, 12).Value = 112513050
$ws.Cells.Item(125, 13).Value = -218040
$ws.Cells.Item(125, 14).Value = -112517970

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 8097.9614
$ws.Cells.Item(61, 9).Value = 8261.879999999999
$ws.Cells.Item(61, 11).Value = 8261.879999999999
$ws.Cells.Item(61, 13).Value = -8049.879999999999
$ws.Cells.Item(63, 8).Value = 670157.75
$ws.Cells.Item(63, 9).Value = 3173.7
$ws.Cells.Item(63, 11).Value = 3173.7
$ws.Cells.Item(63, 13).Value = -2487.7
$ws.Cells.Item(66, 8).Value = 670157.75
$ws.Cells.Item(66, 9).Value = 3173.7
$ws.Cells.Item(66, 11).Value = 15868.5
$ws.Cells.Item(66, 13).Value = -12436.5
$ws.Cells.Item(74, 8).Value = 4650.294
$ws.Cells.Item(74, 9).Value = 3325.3572
$ws.Cells.Item(74, 11).Value = 3325.3572
$ws.Cells.Item(74, 13).Value = -2451.3572
$ws.Cells.Item(77, 8).Value = 4650.294
$ws.Cells.Item(77, 9).Value = 3325.3572
$ws.Cells.Item(77, 11).Value = 16626.786
$ws.Cells.Item(77, 13).Value = -12258.786
$ws.Cells.Item(102, 8).Value = 5697.5713
$ws.Cells.Item(102, 9).Value = 3177.7
$ws.Cells.Item(102, 11).Value = 3177.7
$ws.Cells.Item(102, 13).Value = -1555.7
$ws.Cells.Item(122, 8).Value = 1698.579
$ws.Cells.Item(122, 9).Value = 1704.0555
$ws.Cells.Item(122, 10).Value = 1600
$ws.Cells.Item(122, 11).Value = 5112.166499999999
$ws.Cells.Item(122, 12).Value = 4800
$ws.Cells.Item(122, 13).Value = -2662.166499999999
$ws.Cells.Item(122, 14).Value = -9700
$ws.Cells.Item(132, 8).Value = 2093.5715
$ws.Cells.Item(132, 9).Value = 2093.5715
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 6280.7145
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -3750.7145
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 8097.9614
$ws.Cells.Item(136, 9).Value = 8261.879999999999
$ws.Cells.Item(136, 11).Value = 24785.64
$ws.Cells.Item(136, 13).Value = -22235.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1993.3
$ws.Cells.Item(86, 9).Value = 1939.5
$ws.Cells.Item(86, 10).Value = 2074
$ws.Cells.Item(86, 11).Value = 1939.5
$ws.Cells.Item(86, 12).Value = 2074
$ws.Cells.Item(86, 13).Value = -816.5
$ws.Cells.Item(86, 14).Value = -4320
$ws.Cells.Item(89, 8).Value = 1993.3
$ws.Cells.Item(89, 9).Value = 1939.5
$ws.Cells.Item(89, 10).Value = 2074
$ws.Cells.Item(89, 11).Value = 9697.5
$ws.Cells.Item(89, 12).Value = 10370
$ws.Cells.Item(89, 13).Value = -4081.5
$ws.Cells.Item(89, 14).Value = -21602
$ws.Cells.Item(97, 8).Value = 14211.1
$ws.Cells.Item(97, 9).Value = 8037.8335
$ws.Cells.Item(97, 11).Value = 8037.8335
$ws.Cells.Item(97, 13).Value = -7046.8335
$ws.Cells.Item(99, 8).Value = 5140.25
$ws.Cells.Item(99, 9).Value = 4128.9287
$ws.Cells.Item(99, 11).Value = 4128.9287
$ws.Cells.Item(99, 13).Value = -2630.9287
$ws.Cells.Item(105, 8).Value = 2497.8333
$ws.Cells.Item(105, 9).Value = 2419.7
$ws.Cells.Item(105, 11).Value = 2419.7
$ws.Cells.Item(105, 13).Value = -672.6999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 351.8
$ws.Cells.Item(22, 9).Value = 396
$ws.Cells.Item(22, 11).Value = 396
$ws.Cells.Item(22, 13).Value = -46
$ws.Cells.Item(58, 8).Value = 10106.137
$ws.Cells.Item(58, 10).Value = 13499.77
$ws.Cells.Item(58, 12).Value = 13499.77
$ws.Cells.Item(58, 14).Value = -13905.77
$ws.Cells.Item(74, 8).Value = 41032.25
$ws.Cells.Item(74, 10).Value = 41032.25
$ws.Cells.Item(74, 12).Value = 41032.25
$ws.Cells.Item(74, 14).Value = -42780.25
$ws.Cells.Item(77, 8).Value = 41032.25
$ws.Cells.Item(77, 10).Value = 41032.25
$ws.Cells.Item(77, 12).Value = 123096.75
$ws.Cells.Item(77, 14).Value = -131832.75
$ws.Cells.Item(110, 8).Value = 72097.25
$ws.Cells.Item(110, 10).Value = 79463
$ws.Cells.Item(110, 12).Value = 79463
$ws.Cells.Item(110, 14).Value = -87643
$ws.Cells.Item(112, 8).Value = 46302.145
$ws.Cells.Item(112, 10).Value = 46302.145
$ws.Cells.Item(112, 12).Value = 46302.145
$ws.Cells.Item(112, 14).Value = -49256.145
$ws.Cells.Item(122, 8).Value = 2411.7144
$ws.Cells.Item(122, 9).Value = 2423.2593
$ws.Cells.Item(122, 11).Value = 7269.777900000001
$ws.Cells.Item(122, 13).Value = -4819.777900000001
$ws.Cells.Item(132, 8).Value = 5124.6
$ws.Cells.Item(132, 9).Value = 3644.4443
$ws.Cells.Item(132, 11).Value = 10933.3329
$ws.Cells.Item(132, 13).Value = -8403.332900000001
$ws.Cells.Item(136, 8).Value = 10106.137
$ws.Cells.Item(136, 10).Value = 13499.77
$ws.Cells.Item(136, 12).Value = 40499.31
$ws.Cells.Item(136, 14).Value = -45599.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(128, 8).Value = 199999
$ws.Cells.Item(128, 9).Value = 199999
$ws.Cells.Item(128, 11).Value = 599997
$ws.Cells.Item(128, 13).Value = -595017

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(25, 8).Value = 23925
$ws.Cells.Item(25, 10).Value = 23925
$ws.Cells.Item(25, 12).Value = 23925
$ws.Cells.Item(25, 14).Value = -24983

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3973.7144
$ws.Cells.Item(7, 9).Value = 4023.5789
$ws.Cells.Item(7, 10).Value = 3500
$ws.Cells.Item(7, 11).Value = 4023.5789
$ws.Cells.Item(7, 12).Value = 3500
$ws.Cells.Item(7, 13).Value = -3911.5789
$ws.Cells.Item(7, 14).Value = -3724
$ws.Cells.Item(38, 8).Value = 517500
$ws.Cells.Item(38, 10).Value = 35000
$ws.Cells.Item(38, 12).Value = 35000
$ws.Cells.Item(38, 14).Value = -35820
$ws.Cells.Item(55, 8).Value = 238.25
$ws.Cells.Item(55, 9).Value = 225.41667
$ws.Cells.Item(55, 10).Value = 257.5
$ws.Cells.Item(55, 11).Value = 225.41667
$ws.Cells.Item(55, 12).Value = 257.5
$ws.Cells.Item(55, 13).Value = -52.41667000000001
$ws.Cells.Item(55, 14).Value = -603.5
$ws.Cells.Item(68, 8).Value = 2555.3635
$ws.Cells.Item(68, 9).Value = 1967.8334
$ws.Cells.Item(68, 10).Value = 3260.4
$ws.Cells.Item(68, 11).Value = 1967.8334
$ws.Cells.Item(68, 12).Value = 3260.4
$ws.Cells.Item(68, 13).Value = -1218.8334
$ws.Cells.Item(68, 14).Value = -4758.4
$ws.Cells.Item(71, 8).Value = 2555.3635
$ws.Cells.Item(71, 9).Value = 1967.8334
$ws.Cells.Item(71, 10).Value = 3260.4
$ws.Cells.Item(71, 11).Value = 9839.166999999999
$ws.Cells.Item(71, 12).Value = 16302
$ws.Cells.Item(71, 13).Value = -6095.166999999999
$ws.Cells.Item(71, 14).Value = -23790
$ws.Cells.Item(94, 8).Value = 48997.5
$ws.Cells.Item(94, 10).Value = 48997.5
$ws.Cells.Item(94, 12).Value = 48997.5
$ws.Cells.Item(94, 14).Value = -50349.5
$ws.Cells.Item(126, 8).Value = 3973.7144
$ws.Cells.Item(126, 9).Value = 4023.5789
$ws.Cells.Item(126, 10).Value = 3500
$ws.Cells.Item(126, 11).Value = 12070.7367
$ws.Cells.Item(126, 12).Value = 10500
$ws.Cells.Item(126, 13).Value = -9600.736699999999
$ws.Cells.Item(126, 14).Value = -15440
$ws.Cells.Item(136, 8).Value = 8239.41
$ws.Cells.Item(136, 9).Value = 8238.892
$ws.Cells.Item(136, 11).Value = 24716.676
$ws.Cells.Item(136, 13).Value = -22166.676

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 59999
$ws.Cells.Item(69, 10).Value = 59999
$ws.Cells.Item(69, 12).Value = 59999
$ws.Cells.Item(69, 14).Value = -61497
$ws.Cells.Item(72, 8).Value = 59999
$ws.Cells.Item(72, 10).Value = 59999
$ws.Cells.Item(72, 12).Value = 179997
$ws.Cells.Item(72, 14).Value = -187485
$ws.Cells.Item(100, 8).Value = 661.8823
$ws.Cells.Item(100, 9).Value = 571.2727
$ws.Cells.Item(100, 10).Value = 828
$ws.Cells.Item(100, 11).Value = 1142.5454
$ws.Cells.Item(100, 12).Value = 1656
$ws.Cells.Item(100, 13).Value = -601.5454
$ws.Cells.Item(100, 14).Value = -2738
$ws.Cells.Item(107, 8).Value = 794.2857
$ws.Cells.Item(107, 9).Value = 550.8889
$ws.Cells.Item(107, 11).Value = 1652.6667
$ws.Cells.Item(107, 13).Value = 267.3332999999998
$ws.Cells.Item(122, 8).Value = 4960.2856
$ws.Cells.Item(122, 9).Value = 1996.6666
$ws.Cells.Item(122, 11).Value = 5989.9998
$ws.Cells.Item(122, 13).Value = -3539.9998
$ws.Cells.Item(132, 8).Value = 4127.619
$ws.Cells.Item(132, 10).Value = 7082.778
$ws.Cells.Item(132, 12).Value = 21248.334
$ws.Cells.Item(132, 14).Value = -26308.334
$ws.Cells.Item(135, 8).Value = 98995
$ws.Cells.Item(135, 10).Value = 98995
$ws.Cells.Item(135, 12).Value = 98995
$ws.Cells.Item(135, 14).Value = -109135
$ws.Cells.Item(136, 8).Value = 2431.6
$ws.Cells.Item(136, 9).Value = 1875.7059
$ws.Cells.Item(136, 11).Value = 5627.1177
$ws.Cells.Item(136, 13).Value = -3077.1177
